$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style as the other header cells (e.g. H1) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data for rows 2-30: column I, column J
$values = @(
    @(5, 6),
    @(1, 2),
    @(1, 4),
    @(10, 10),
    @(1, 3),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(3, 7),
    @(5, 7),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(4, 6),
    @(5, 6),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(1, 2),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
